$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 3364.8333
$ws.Range("I11").Value = 3364.8333
$ws.Range("K11").Value = 3364.8333
$ws.Range("M11").Value = -3224.8333
$ws.Range("H42").Value = 192.73334
$ws.Range("I42").Value = 60.625
$ws.Range("J42").Value = 343.7143
$ws.Range("K42").Value = 181.875
$ws.Range("L42").Value = 1031.1429
$ws.Range("M42").Value = 48.125
$ws.Range("N42").Value = -1491.1429
$ws.Range("H107").Value = 19608544
$ws.Range("I107").Value = 22222880
$ws.Range("K107").Value = 22222880
$ws.Range("M107").Value = -22220960
$ws.Range("H132").Value = 47622900
$ws.Range("I132").Value = 55559550
$ws.Range("K132").Value = 166678650
$ws.Range("M132").Value = -166676120
$ws.Range("H137").Value = 151435.42
$ws.Range("I137").Value = 298079.16
$ws.Range("K137").Value = 894237.48
$ws.Range("M137").Value = -891687.48
$ws.Range("H138").Value = 2445.0308
$ws.Range("J138").Value = 4414.407
$ws.Range("L138").Value = 13243.221
$ws.Range("N138").Value = -23523.221
$ws.Range("H141").Value = 2360.5881
$ws.Range("I141").Value = 2043.6364
$ws.Range("J141").Value = 2941.6667
$ws.Range("K141").Value = 6130.9092
$ws.Range("L141").Value = 8825.000100000001
$ws.Range("M141").Value = -950.9092000000001
$ws.Range("N141").Value = -19185.0001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H10").Value = 7500
$ws.Range("I10").Value = 7500
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 7500
$ws.Range("L10").Value = 0
$ws.Range("M10").Value = -7330
$ws.Range("N10").Value = ""
$ws.Range("H61").Value = 2120.7104
$ws.Range("I61").Value = 1934.2916
$ws.Range("J61").Value = 2440.2856
$ws.Range("K61").Value = 1934.2916
$ws.Range("L61").Value = 2440.2856
$ws.Range("M61").Value = -1722.2916
$ws.Range("N61").Value = -2864.2856
$ws.Range("H74").Value = 41998.25
$ws.Range("I74").Value = 4293.2104
$ws.Range("K74").Value = 4293.2104
$ws.Range("M74").Value = -3419.2104
$ws.Range("H77").Value = 41998.25
$ws.Range("I77").Value = 4293.2104
$ws.Range("K77").Value = 21466.052
$ws.Range("M77").Value = -17098.052
$ws.Range("H132").Value = 2193.5715
$ws.Range("I132").Value = 2028.4546
$ws.Range("J132").Value = 2473
$ws.Range("K132").Value = 6085.3638
$ws.Range("L132").Value = 7419
$ws.Range("M132").Value = -3555.3638
$ws.Range("N132").Value = -12479
$ws.Range("H136").Value = 2120.7104
$ws.Range("I136").Value = 1934.2916
$ws.Range("J136").Value = 2440.2856
$ws.Range("K136").Value = 5802.8748
$ws.Range("L136").Value = 7320.8568
$ws.Range("M136").Value = -3252.8748
$ws.Range("N136").Value = -12420.8568

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 1489.1666
$ws.Range("I5").Value = 450
$ws.Range("J5").Value = 2008.75
$ws.Range("K5").Value = 450
$ws.Range("L5").Value = 2008.75
$ws.Range("M5").Value = -337
$ws.Range("N5").Value = -2234.75
$ws.Range("H24").Value = 1151.25
$ws.Range("I24").Value = 1151.25
$ws.Range("J24").Value = 0
$ws.Range("K24").Value = 1151.25
$ws.Range("L24").Value = 0
$ws.Range("M24").Value = -916.25
$ws.Range("N24").Value = ""
$ws.Range("H86").Value = 5889792
$ws.Range("J86").Value = 1755
$ws.Range("L86").Value = 1755
$ws.Range("N86").Value = -4001
$ws.Range("H89").Value = 5889792
$ws.Range("J89").Value = 1755
$ws.Range("L89").Value = 8775
$ws.Range("N89").Value = -20007
$ws.Range("H134").Value = 3509.5527
$ws.Range("I134").Value = 1667.862
$ws.Range("K134").Value = 5003.586
$ws.Range("M134").Value = -2468.586

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2176.4
$ws.Range("I16").Value = 1721.75
$ws.Range("K16").Value = 1721.75
$ws.Range("M16").Value = -1434.75
$ws.Range("H107").Value = 1238.9117
$ws.Range("I107").Value = 1149.0385
$ws.Range("J107").Value = 1531
$ws.Range("K107").Value = 1149.0385
$ws.Range("L107").Value = 1531
$ws.Range("M107").Value = 770.9614999999999
$ws.Range("N107").Value = -5371
$ws.Range("H113").Value = 2176.4
$ws.Range("I113").Value = 1721.75
$ws.Range("K113").Value = 1721.75
$ws.Range("M113").Value = 448.25
$ws.Range("H122").Value = 2959.818
$ws.Range("I122").Value = 2601.1428
$ws.Range("J122").Value = 3587.5
$ws.Range("K122").Value = 7803.428400000001
$ws.Range("L122").Value = 10762.5
$ws.Range("M122").Value = -5353.428400000001
$ws.Range("N122").Value = -15662.5
$ws.Range("H132").Value = 37831.48
$ws.Range("I132").Value = 2299.9565
$ws.Range("K132").Value = 6899.869499999999
$ws.Range("M132").Value = -4369.869499999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 46636.75
$ws.Range("J37").Value = 46636.75
$ws.Range("L37").Value = 139910.25
$ws.Range("N37").Value = -140134.25
$ws.Range("H129").Value = 2222982
$ws.Range("I129").Value = 2857776.8
$ws.Range("J129").Value = 1200
$ws.Range("K129").Value = 8573330.399999999
$ws.Range("L129").Value = 3600
$ws.Range("M129").Value = -8568330.399999999
$ws.Range("N129").Value = -13600
$ws.Range("H133").Value = 0
$ws.Range("I133").Value = 0
$ws.Range("K133").Value = 0
$ws.Range("M133").Value = ""
$ws.Range("H134").Value = 3764.0715
$ws.Range("J134").Value = 15000
$ws.Range("L134").Value = 45000
$ws.Range("N134").Value = -55140
$ws.Range("H139").Value = 1788
$ws.Range("I139").Value = 1329.1428
$ws.Range("K139").Value = 3987.4284
$ws.Range("M139").Value = 1152.5716
$ws.Range("H140").Value = 1610.1666
$ws.Range("I140").Value = 1392.9546
$ws.Range("K140").Value = 4178.8638
$ws.Range("M140").Value = 1001.1362

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 6781532
$ws.Range("I126").Value = 3033518.2
$ws.Range("J126").Value = 20836582
$ws.Range("K126").Value = 9100554.600000001
$ws.Range("L126").Value = 62509746
$ws.Range("M126").Value = -9098084.600000001
$ws.Range("N126").Value = -62514686
$ws.Range("H132").Value = 4734.3335
$ws.Range("I132").Value = 3797.4
$ws.Range("J132").Value = 5905.5
$ws.Range("K132").Value = 11392.2
$ws.Range("L132").Value = 17716.5
$ws.Range("M132").Value = -8862.200000000001
$ws.Range("N132").Value = -22776.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 34545.906
$ws.Range("I136").Value = 44582.125
$ws.Range("J136").Value = 4437.25
$ws.Range("K136").Value = 133746.375
$ws.Range("L136").Value = 13311.75
$ws.Range("M136").Value = -131196.375
$ws.Range("N136").Value = -18411.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H3").Value = 7500000
$ws.Range("I3").Value = 7500000
$ws.Range("K3").Value = 7500000
$ws.Range("M3").Value = -7499886
$ws.Range("H45").Value = 37998.5
$ws.Range("J45").Value = 37998.5
$ws.Range("L45").Value = 37998.5
$ws.Range("N45").Value = -38980.5
$ws.Range("H122").Value = 2257.6924
$ws.Range("I122").Value = 1561.375
$ws.Range("J122").Value = 3371.8
$ws.Range("K122").Value = 4684.125
$ws.Range("L122").Value = 10115.4
$ws.Range("M122").Value = -2234.125
$ws.Range("N122").Value = -15015.4
$ws.Range("H132").Value = 41710200
$ws.Range("I132").Value = 50007276
$ws.Range("K132").Value = 150021828
$ws.Range("M132").Value = -150019298

Write-Host "Applied 36 row updates across 8 sheets"
